$wb = $excel.ActiveWorkbook

$wsCoord = $wb.Worksheets.Item("coordinates")
$wsSource = $wb.Worksheets.Item("source")

# Update A/B values in "coordinates" sheet
$wsCoord.Cells.Item(2, 2).Value = 218

$wsCoord.Cells.Item(3, 1).Value = 347
$wsCoord.Cells.Item(3, 2).Value = 153

$wsCoord.Cells.Item(4, 1).Value = 347
$wsCoord.Cells.Item(4, 2).Value = 303

$wsCoord.Cells.Item(5, 1).Value = 347
$wsCoord.Cells.Item(5, 2).Value = 453

$wsCoord.Cells.Item(6, 1).Value = 347
$wsCoord.Cells.Item(6, 2).Value = 603

$wsCoord.Cells.Item(7, 1).Value = 498
$wsCoord.Cells.Item(7, 2).Value = 228

$wsCoord.Cells.Item(8, 1).Value = 498
$wsCoord.Cells.Item(8, 2).Value = 378

$wsCoord.Cells.Item(9, 1).Value = 498
$wsCoord.Cells.Item(9, 2).Value = 528

$wsCoord.Cells.Item(10, 1).Value = 649
$wsCoord.Cells.Item(10, 2).Value = 153

$wsCoord.Cells.Item(11, 1).Value = 649
$wsCoord.Cells.Item(11, 2).Value = 303

$wsCoord.Cells.Item(12, 1).Value = 649
$wsCoord.Cells.Item(12, 2).Value = 453

$wsCoord.Cells.Item(13, 1).Value = 800
$wsCoord.Cells.Item(13, 2).Value = 78

$wsCoord.Cells.Item(14, 1).Value = 800
$wsCoord.Cells.Item(14, 2).Value = 228

$wsCoord.Cells.Item(15, 1).Value = 800
$wsCoord.Cells.Item(15, 2).Value = 378

$wsCoord.Cells.Item(16, 1).Value = 951
$wsCoord.Cells.Item(16, 2).Value = 153

$wsCoord.Cells.Item(17, 1).Value = 951
$wsCoord.Cells.Item(17, 2).Value = 303

$wsCoord.Cells.Item(18, 1).Value = 951
$wsCoord.Cells.Item(18, 2).Value = 453

$wsCoord.Cells.Item(19, 1).Value = 1102
$wsCoord.Cells.Item(19, 2).Value = 228

$wsCoord.Cells.Item(20, 1).Value = 1102
$wsCoord.Cells.Item(20, 2).Value = 378

$wsCoord.Cells.Item(21, 1).Value = 1102
$wsCoord.Cells.Item(21, 2).Value = 528

$wsCoord.Cells.Item(22, 1).Value = 1253
$wsCoord.Cells.Item(22, 2).Value = 303

$wsCoord.Cells.Item(23, 1).Value = 1253
$wsCoord.Cells.Item(23, 2).Value = 453

$wsCoord.Cells.Item(24, 1).Value = 1404
$wsCoord.Cells.Item(24, 2).Value = 228

# Update "source" sheet value
$wsSource.Cells.Item(2, 2).Value = 21

# Sheet view / selection changes
$wsCoord.Range("K15").Select() | Out-Null
$wsSource.Range("C4").Select() | Out-Null

# Make "source" the active sheet (tab selected) last so it becomes the active tab
$wsSource.Activate() | Out-Null
